# Update cryptos list data (prices & volume changes) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.127.83"
$ws.Range("E2").Value = "'  -0.41%  "
$ws.Range("D3").Value = "'1.778.90"
$ws.Range("E3").Value = "'  -2.54%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "'  +0.19%  "
$ws.Range("D5").Value = "'225.38"
$ws.Range("E5").Value = "'  -1.67%  "
$ws.Range("E6").Value = "'  +0.43%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "'  +0.10%  "
$ws.Range("D8").Value = "'31.65"
$ws.Range("E8").Value = "'  +0.22%  "
$ws.Range("E9").Value = "'  -0.62%  "
$ws.Range("E10").Value = "'  -1.94%  "
$ws.Range("E11").Value = "'  -0.09%  "
$ws.Range("D12").Value = "'2.034.39"
$ws.Range("E12").Value = "'  -2.33%  "
$ws.Range("E13").Value = "'  +6.19%  "
$ws.Range("D14").Value = "'1.779.26"
$ws.Range("E14").Value = "'  -2.64%  "
$ws.Range("D15").Value = "'0.624"
$ws.Range("E15").Value = "'  -3.25%  "
$ws.Range("D16").Value = "'34.142.12"
$ws.Range("E16").Value = "'  -0.10%  "
$ws.Range("E17").Value = "'  -1.70%  "
$ws.Range("D18").Value = "'68.77"
$ws.Range("E18").Value = "'  -1.41%  "
$ws.Range("D19").Value = "'254.02"
$ws.Range("E19").Value = "'  -2.02%  "
$ws.Range("D20").Value = "'0.0₃0737"
$ws.Range("E20").Value = "'  -2.02%  "
$ws.Range("E21").Value = "'  -0.03%  "
$ws.Range("D22").Value = "'10.35"
$ws.Range("E22").Value = "'  -2.41%  "
$ws.Range("E23").Value = "'  -3.79%  "
$ws.Range("E24").Value = "'  -3.90%  "
$ws.Range("D25").Value = "'156.69"
$ws.Range("E25").Value = "'  -1.07%  "
$ws.Range("D26").Value = "'16.42"
$ws.Range("E26").Value = "'  -1.48%  "
$ws.Range("D27").Value = "'6.98"
$ws.Range("E27").Value = "'  -2.56%  "
$ws.Range("E28").Value = "'  -1.39%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "'  +0.09%  "
$ws.Range("E30").Value = "'  -3.05%  "
$ws.Range("E31").Value = "'  -0.63%  "
$ws.Range("E32").Value = "'  -1.28%  "
$ws.Range("D33").Value = "'3.57"
$ws.Range("E33").Value = "'  +0.49%  "
$ws.Range("E34").Value = "'  +1.61%  "
$ws.Range("D35").Value = "'1.439.75"
$ws.Range("E35").Value = "'  -7.07%  "
$ws.Range("E36").Value = "'  -3.91%  "
$ws.Range("B37").Value = "'VeChain"
$ws.Range("C37").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.0187"
$ws.Range("E37").Value = "'  -0.34%  "
$ws.Range("B38").Value = "'ImmutableX"
$ws.Range("C38").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'0.625"
$ws.Range("E38").Value = "'  -0.53%  "
$ws.Range("D39").Value = "'2.86"
$ws.Range("E39").Value = "'  +1.22%  "
$ws.Range("D40").Value = "'82.82"
$ws.Range("E40").Value = "'  -3.30%  "
$ws.Range("E41").Value = "'  +0.17%  "
$ws.Range("D42").Value = "'0.889"
$ws.Range("E42").Value = "'  -3.20%  "
$ws.Range("E43").Value = "'  -5.14%  "
$ws.Range("D44").Value = "'0.0509"
$ws.Range("E44").Value = "'  -3.06%  "
$ws.Range("E45").Value = "'  -1.88%  "
$ws.Range("B46").Value = "'RocketPoolETH"
$ws.Range("C46").Value = "'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "'1.934.86"
$ws.Range("E46").Value = "'  -2.46%  "
$ws.Range("B47").Value = "'FraxShare"
$ws.Range("C47").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'5.80"
$ws.Range("E47").Value = "'  +1.12%  "
$ws.Range("D48").Value = "'12.02"
$ws.Range("E48").Value = "'  +1.70%  "
$ws.Range("E49").Value = "'  +0.08%  "
$ws.Range("D50").Value = "'98.28"
$ws.Range("E50").Value = "'  +2.11%  "
$ws.Range("D51").Value = "'49.70"
$ws.Range("E51").Value = "'  -6.38%  "
